# Apply the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.736.97"
$ws.Range("E2").Value = "  -1.08%  "

# Row 3
$ws.Range("D3").Value = "3.514.94"
$ws.Range("E3").Value = "  -2.67%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.44"
$ws.Range("E5").Value = "  +0.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.83"
$ws.Range("E6").Value = "  -2.60%  "

# Row 7
$ws.Range("D7").Value = "3.507.70"
$ws.Range("E7").Value = "  -2.77%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.615"
$ws.Range("E8").Value = "  -2.89%  "

# Row 9
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.202"
$ws.Range("E10").Value = "  +10.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.649"
$ws.Range("E11").Value = "  -2.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.05"
$ws.Range("E12").Value = "  -3.99%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000307"
$ws.Range("E13").Value = "  -1.51%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.49"
$ws.Range("E14").Value = "  -2.65%  "

# Row 15
$ws.Range("D15").Value = "4.074.44"
$ws.Range("E15").Value = "  -2.75%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.35"
$ws.Range("E16").Value = "  -3.70%  "

# Row 17
$ws.Range("D17").Value = "69.650.79"
$ws.Range("E17").Value = "  -1.09%  "

# Row 18
$ws.Range("D18").Value = "3.499.81"
$ws.Range("E18").Value = "  -3.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.34"
$ws.Range("E19").Value = "  -3.47%  "

# Row 20
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "529.70"
$ws.Range("E21").Value = "  +9.44%  "

# Row 22
$ws.Range("E22").Value = "  -3.76%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.34"
$ws.Range("E23").Value = "  -3.97%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.93"
$ws.Range("E24").Value = "  -2.59%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.49"
$ws.Range("E25").Value = "  +1.79%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.32"
$ws.Range("E26").Value = "  -1.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.18"
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.97"
$ws.Range("E28").Value = "  -1.85%  "

# Row 29
$ws.Range("E29").Value = "  -3.80%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.91"
$ws.Range("E30").Value = "  -1.60%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.38"
$ws.Range("E31").Value = "  -3.72%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.62"
$ws.Range("E32").Value = "  +2.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.89"
$ws.Range("E33").Value = "  -2.24%  "

# Row 34
$ws.Range("E34").Value = "  -5.09%  "

# Row 35
$ws.Range("B35").Value = "TheGraph"
$ws.Range("C35").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.419"
$ws.Range("E35").Value = "  +4.97%  "

# Row 36
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "537.26"
$ws.Range("E36").Value = "  -8.26%  "

# Row 37
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.08"
$ws.Range("E37").Value = "  +7.89%  "

# Row 38
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.93"
$ws.Range("E38").Value = "  -3.62%  "

# Row 39
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0766"
$ws.Range("E40").Value = "  -6.65%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.17"
$ws.Range("E41").Value = "  -6.45%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  -2.88%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  -2.60%  "

# Row 44
$ws.Range("D44").Value = "3.333.61"
$ws.Range("E44").Value = "  +2.95%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.99"
$ws.Range("E45").Value = "  -4.00%  "

# Row 46
$ws.Range("E46").Value = "  -1.62%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.46"
$ws.Range("E47").Value = "  +2.22%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.135"
$ws.Range("E48").Value = "  -3.40%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.90"
$ws.Range("E49").Value = "  -7.83%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.48"
$ws.Range("E51").Value = "  +2.44%  "
